# Matriz Compiladores2.xlsx — "MAS Arreglada para igual"
#
# The row for state 8 ("FinSimboloIncPos") had its label corrected to a new,
# distinct token "FinSimbolo" across the whole row (every destination-state
# column except the I column, which legitimately keeps the old
# "FinSimboloIncPos" label). Cell O8 also had a stray/mismatched cell format
# (left over from a copy-paste) that gets normalized to match the rest of the
# row. Finally the view was scrolled/zoomed out and the selection left on the
# edited row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Re-label row 8 (state 8) with the corrected token "FinSimbolo".
#    Column I keeps its original "FinSimboloIncPos" value untouched.
$cols = @("B", "C", "D", "E", "F", "G", "H", "J", "K", "L", "M", "N", "O", "P", "Q", "R")
foreach ($col in $cols) {
    $ws.Range($col + "8").Value2 = "FinSimbolo"
}

# 2. O8 had a slightly different (but visually identical) cell style than its
#    neighbours — normalize it by copying the format from N8 onto O8.
$ws.Range("N8").Copy() | Out-Null
$ws.Range("O8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 3. Zoom the sheet out to 62%.
$excel.ActiveWindow.Zoom = 62

# 4. Leave the selection on the row that was just edited.
$ws.Range("B8:R8").Select() | Out-Null
